$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text: column E ("eight year graduation count" -> "eightYearsGradCount")
# and column H ("eight year graduation rate" -> "eightYearsGradRate").
$ws.Range("E1").Value = "eightYearsGradCount"
$ws.Range("H1").Value = "eightYearsGradRate"

# Move the active selection to H1 (matches the updated sheetView selection).
$ws.Range("H1").Select()
